$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.374.12"
$ws.Range("E2").Value = "  +6.57%  "
$ws.Range("D3").Value = "3.115.97"
$ws.Range("E3").Value = "  +4.42%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "585.61"
$ws.Range("E5").Value = "  +3.36%  "
$ws.Range("E6").Value = "  +4.86%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.107.77"
$ws.Range("E8").Value = "  +4.54%  "
$ws.Range("D9").Value = "0.531"
$ws.Range("E10").Value = "  +13.88%  "
$ws.Range("D11").Value = "5.80"
$ws.Range("E11").Value = "  +7.73%  "
$ws.Range("E12").Value = "  +4.15%  "
$ws.Range("E13").Value = "  +8.24%  "
$ws.Range("E14").Value = "  +5.77%  "
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "3.630.48"
$ws.Range("E16").Value = "  +4.41%  "
$ws.Range("D17").Value = "7.18"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").Value = "63.268.41"
$ws.Range("E18").Value = "  +6.39%  "
$ws.Range("D19").Value = "3.112.56"
$ws.Range("E19").Value = "  +4.42%  "
$ws.Range("D20").Value = "468.31"
$ws.Range("E20").Value = "  +7.21%  "
$ws.Range("D21").Value = "14.16"
$ws.Range("E21").Value = "  +4.19%  "
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("E23").Value = "  +7.57%  "
$ws.Range("D24").Value = "13.33"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "82.12"
$ws.Range("E25").Value = "  +2.68%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +10.52%  "
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("E29").Value = "  +5.24%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  +10.91%  "
$ws.Range("D32").Value = "26.98"
$ws.Range("E32").Value = "  +4.81%  "
$ws.Range("E33").Value = "  +4.47%  "
$ws.Range("D34").Value = "0.0₃0867"
$ws.Range("E34").Value = "  +12.08%  "
$ws.Range("E35").Value = "  +16.48%  "
$ws.Range("E36").Value = "  +6.78%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "3.34"
$ws.Range("E37").Value = "  +20.02%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "6.09"
$ws.Range("E38").Value = "  +3.27%  "
$ws.Range("E39").Value = "  +4.05%  "
$ws.Range("D40").Value = "440.50"
$ws.Range("E40").Value = "  +9.94%  "
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").Value = "2.927.71"
$ws.Range("E42").Value = "  +7.03%  "
$ws.Range("E43").Value = "  +5.39%  "
$ws.Range("E44").Value = "  +11.97%  "
$ws.Range("E45").Value = "  +5.02%  "
$ws.Range("E46").Value = "  +8.53%  "
$ws.Range("D47").Value = "35.55"
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").Value = "123.26"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Value = "24.72"
$ws.Range("E51").Value = "  +6.13%  "
